# "stateless entities outside the US"
#
# A new metric column pair, "IMF (20%)", is introduced ahead of the existing
# "IMF" column pair (and removes the previously-present "OECD (20%)" pair),
# for each of the three eight-column blocks (Sales/M_%cit, M_ETR, M_PL).
# Concretely, per block the columns read:
#   ... GFA - Sales | GFA - Sales + Emp | IMF - Sales | IMF - Sales + Emp | OECD (20%) - Sales | OECD (20%) - Sales + Emp | OECD - Sales | OECD - Sales + Emp
# becomes
#   ... GFA - Sales | GFA - Sales + Emp | IMF (20%) - Sales | IMF (20%) - Sales + Emp | IMF - Sales | IMF - Sales + Emp | OECD - Sales | OECD - Sales + Emp
#
# The old "IMF" figures slide over into the (old) "OECD (20%)" slot unchanged,
# while the vacated "IMF" slot is populated with the new "IMF (20%)" figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header relabeling (row 2, one pair of columns per eight-column block) ---
foreach ($addr in @("D2", "L2", "T2")) {
    $ws.Range($addr).Value = "IMF (20%) - Sales"
}
foreach ($addr in @("E2", "M2", "U2")) {
    $ws.Range($addr).Value = "IMF (20%) - Sales + Emp"
}
foreach ($addr in @("F2", "N2", "V2")) {
    $ws.Range($addr).Value = "IMF - Sales"
}
foreach ($addr in @("G2", "O2", "W2")) {
    $ws.Range($addr).Value = "IMF - Sales + Emp"
}

# --- Data: block 1 (M_%cit, columns B-I) ---
$ws.Range("D4").Value = 0.1442535571439624
$ws.Range("E4").Value = 0.2862580632476772
$ws.Range("F4").Value = 0.7212677857198116
$ws.Range("G4").Value = 1.431290316238386

$ws.Range("D5").Value = 0.6761226210377156
$ws.Range("E5").Value = 0.7404358172653206
$ws.Range("F5").Value = 3.380613105188575
$ws.Range("G5").Value = 3.702179086326601

$ws.Range("D6").Value = 1.054230044300623
$ws.Range("E6").Value = 2.123036393002244
$ws.Range("F6").Value = 5.271150221503113
$ws.Range("G6").Value = 10.61518196501122

$ws.Range("D7").Value = 0.2076904293943114
$ws.Range("E7").Value = 0.4662097981339481
$ws.Range("F7").Value = 1.038452146971557
$ws.Range("G7").Value = 2.331048990669741

$ws.Range("D8").Value = 1.508446806846264
$ws.Range("E8").Value = 1.052810667356595
$ws.Range("F8").Value = 7.542234034231321
$ws.Range("G8").Value = 5.264053336782998

$ws.Range("E9").Value = 74.30184066874688
$ws.Range("G9").Value = 371.5092033437344
$ws.Range("F9").ClearContents()

$ws.Range("D10").Value = 0.8575958516640255
$ws.Range("E10").Value = 1.051093146802996
$ws.Range("F10").Value = 4.287979258320128
$ws.Range("G10").Value = 5.25546573401498

# --- Data: block 2 (M_ETR, columns J-Q) ---
$ws.Range("N4").Value = 0.3248011131580942
$ws.Range("O4").Value = 0.2982077124743382

$ws.Range("N5").Value = 0.2960286019796783
$ws.Range("O5").Value = 0.2955367085519889

$ws.Range("N6").Value = 0.6442414585966483
$ws.Range("O6").Value = 0.6271429106953286

$ws.Range("N7").Value = 0.6260699716191001
$ws.Range("O7").Value = 0.6260699716191001

$ws.Range("O9").Value = 0.33845850318029
$ws.Range("N9").ClearContents()

$ws.Range("N10").Value = 0.9678569072646386
$ws.Range("O10").Value = 0.9487490577047336

# --- Data: block 3 (M_PL, columns R-Y) ---
$ws.Range("V4").Value = 30235087160
$ws.Range("W4").Value = 37431807140

$ws.Range("V5").Value = 66666309306
$ws.Range("W5").Value = 66917200922

$ws.Range("V6").Value = 15485301894
$ws.Range("W6").Value = 16245802182

$ws.Range("V7").Value = 1403512928
$ws.Range("W7").Value = 1403512928

$ws.Range("W9").Value = 12096909667
$ws.Range("V9").ClearContents()

$ws.Range("V10").Value = 1986738567
$ws.Range("W10").Value = 2032015673
